# Apply updated dSF (column F) values for the specified rows.
# Map of row number -> new value for column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    6  = -4
    9  = 1
    11 = -4
    12 = -4
    16 = -3
    20 = -2
    23 = -3
    26 = -2
    31 = -2
    39 = -2
    42 = -1
    47 = -3
    50 = 0
    54 = -1
    67 = -13
    68 = -4
    71 = -11
    72 = -1
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
